$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "35.077.11"
$ws.Range("E2").Value = "  +1.14%  "
$ws.Range("D3").Value = "1.848.78"
$ws.Range("E3").Value = "  +2.13%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "233.54"
$ws.Range("E5").Value = "  +0.47%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.622"
$ws.Range("E6").Value = "  +2.93%  "
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "42.29"
$ws.Range("E8").Value = "  +7.47%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.330"
$ws.Range("E9").Value = "  +1.21%  "
$ws.Range("E10").Value = "  +2.10%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("D12").Value = "2.116.08"
$ws.Range("E12").Value = "  +2.12%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.49"
$ws.Range("E13").Value = "  +3.01%  "
$ws.Range("D14").Value = "1.852.35"
$ws.Range("E14").Value = "  +2.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.677"
$ws.Range("E15").Value = "  +0.44%  "
$ws.Range("D17").Value = "35.107.20"
$ws.Range("E17").Value = "  +1.20%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "70.03"
$ws.Range("E18").Value = "  +0.87%  "
$ws.Range("D19").Value = "0.0₃0795"
$ws.Range("E19").Value = "  +1.32%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "240.92"
$ws.Range("E20").Value = "  +0.27%  "
$ws.Range("E21").Value = "  +1.93%  "
$ws.Range("E22").Value = "  +1.93%  "
$ws.Range("E23").Value = "  +0.08%  "
$ws.Range("E24").Value = "  +3.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.05"
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.96"
$ws.Range("E26").Value = "  +3.29%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.58"
$ws.Range("E27").Value = "  +2.17%  "
$ws.Range("E28").Value = "  +3.90%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.72"
$ws.Range("E29").Value = "  +11.22%  "
$ws.Range("E30").Value = "  +2.11%  "
$ws.Range("E31").Value = "  -0.07%  "
$ws.Range("E32").Value = "  -1.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.96"
$ws.Range("E33").Value = "  +0.01%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.65"
$ws.Range("E34").Value = "  +24.75%  "
$ws.Range("E35").Value = "  +10.53%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.767"
$ws.Range("E36").Value = "  +9.91%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.23"
$ws.Range("E37").Value = "  -2.57%  "
$ws.Range("E38").Value = "  +11.48%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "90.79"
$ws.Range("E39").Value = "  -0.88%  "
$ws.Range("E40").Value = "  +4.59%  "
$ws.Range("D41").Value = "1.347.71"
$ws.Range("E41").Value = "  +1.81%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "14.70"
$ws.Range("E42").Value = "  +4.07%  "
$ws.Range("B43").Value = "Gas"
$ws.Range("C43").Value = "https://coinranking.com/coin/hfw0nnnLtSFc7+gas-gas"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "12.91"
$ws.Range("E43").Value = "  +86.37%  "
$ws.Range("B44").Value = "RenderToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.32"
$ws.Range("E44").Value = "  +4.66%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.39"
$ws.Range("E45").Value = "  -3.27%  "
$ws.Range("B46").Value = "MXToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.76"
$ws.Range("E46").Value = "  +3.51%  "
$ws.Range("B47").Value = "FraxShare"
$ws.Range("C47").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "6.37"
$ws.Range("E47").Value = "  +1.44%  "
$ws.Range("B48").Value = "Kaspa"
$ws.Range("C48").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0531"
$ws.Range("E48").Value = "  +3.68%  "
$ws.Range("D49").Value = "2.030.07"
$ws.Range("E49").Value = "  +1.57%  "
$ws.Range("E50").Value = "  +16.32%  "
$ws.Range("E51").Value = "  +0.97%  "
